$wb = $excel.ActiveWorkbook

# 1. Insert a new row into the 'storage_medium' lookup sheet for the new
#    "Nuclease-free water" option, right after "CMC" (row 10) and before
#    "2% PFA/2.5% Glutaraldehyde" (old row 11).
$smws = $wb.Worksheets.Item("storage_medium")
$smws.Range("A11:B11").Insert()
$smws.Range("A11").Value = "Nuclease-free water"
$smws.Range("B11").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000482"

# 2. Update the data validation on the main sheet's storage_medium column (M)
#    so its source range covers the new row (was $A$1:$A$25, now $A$1:$A$26).
$mainws = $wb.Worksheets.Item("Sample Suspension")
$mainws.Range("M2:M1001").Validation.Formula1 = "'storage_medium'!`$A`$1:`$A`$26"

# 3. Bump the .metadata sheet's pav:createdOn timestamp to reflect the edit.
$metaws = $wb.Worksheets.Item(".metadata")
$metaws.Cells.Item(2, 3).Value = "2026-02-24T15:31:04-08:00"
